# "maj template comment à la fin" — move the Comment column (header + its
# row2 comment + row3 type + row4 format hint) from column J to the end of
# the block (column P), shifting SSC/FSC/ParticlesFilter/Labelling/
# AntibodyName/ResultCellCount one column to the left (J..O) to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers: SSC, FSC, ParticlesFilter, Labelling, AntibodyName,
# ResultCellCount move left into J..O; Comment becomes the last header (P).
$ws.Range("J1").Value = "SSC"
$ws.Range("K1").Value = "FSC"
$ws.Range("L1").Value = "ParticlesFilter"
$ws.Range("M1").Value = "Labelling"
$ws.Range("N1").Value = "AntibodyName"
$ws.Range("O1").Value = "ResultCellCount"
$ws.Range("P1").Value = "Comment"

# Row 2 - per-column comment: the "# Commentaire" note moves from J2 to P2;
# the columns it vacates (J..O) become blank like the rest of that block.
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = "# Commentaire"

# Row 3 - type marker: the block stays "#string" except the new last column
# (now ResultCellCount, O3) flips to "#integer" and the vacated Comment slot
# (now P3, ex-ResultCellCount's "#integer") becomes "#string".
$ws.Range("J3").Value = "#string"
$ws.Range("K3").Value = "#string"
$ws.Range("L3").Value = "#string"
$ws.Range("M3").Value = "#string"
$ws.Range("N3").Value = "#string"
$ws.Range("O3").Value = "#integer"
$ws.Range("P3").Value = "#string"

# Row 4 - format hint: "# format: texte libre" moves from J4 to P4; the
# columns it vacates (J..O) become blank.
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = "# format: texte libre"

# Row 5 (example row) was already blank across J..P and stays blank.
